# Applies the Simplified-Chinese translation pass described by the commit
# "New translations removal of usdt tether omni - reminder email to
# clients.docx (Chinese Simplified)" to the Traditional-Chinese draft.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

Replace-Text "提醒 ROW 客戶的電子郵件" "针对ROW客户的提醒邮件"
Replace-Text "主題：" "主题:"
Replace-Text "我們將於 9 月 29 日移除 Tether Omni (USDT)" "我们将在9月29日移除Tether Omni（USDT）"
Replace-Text "向 Tether Omni 道別" "向Tether Omni说再见"
Replace-Text "我們將停止在 Deriv 上提供 Tether Omni (USDT) 作為帳戶貨幣，自 2023/09/29 (00:00 GMT) 起生效。 這是因為 Tether 停止支持 USDT 轉帳的 Omni。" "我们将自2023年9月29日（格林威治标准时间00:00）起，停止在Deriv账户中提供Tether Omni（USDT）作为账户货币。 这是因为Tether已停止支持USDT的Omni转账。"
Replace-Text "我需要做什麼？" "我需要做什么？"
Replace-Text "如果您的 USDT 帳戶有餘額 " "如果您在USDT账户"
Replace-Text "[帳戶 ID]" "[账户ID]"
Replace-Text "，請在上述日期之前提取您的餘額。 如果您有未平倉的頭寸，請在提取餘額之前先關閉它們。" "中有余额，请在上述日期之前提取余额。 如果您有未平仓头寸，请在提取余额之前先关闭它们。"
Replace-Text "檢查我的帳戶" "查看我的账户"
Replace-Text "您的 USDT 帳戶將於 2023/09/29 00:00 GMT 關閉。 任何未平倉頭寸將自動關閉，帳戶餘額將在上述日期後轉移到您最後一個有效的帳戶中" "您的USDT账户将在2023年9月29日00:00 GMT关闭。 任何未平仓头寸将在提到的日期后自动关闭，账户余额将转移到您最后活跃的账户"
Replace-Text "在此過程中將適用標準匯率和費用。" "在此过程中将适用标准汇率和费用。"
Replace-Text "如有任何疑問，請聯繫我們：" "如有任何疑问，请联系我们："
Replace-Text "即時聊天" "实时聊天"

Write-Output "done"
